# New API Query - 2023 Included
# API query to UN performed 11/26/2023.
# Query modified to include 2023 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# short-url column (B) changed for the data rows 2-6
$ws.Range("B2:B6").Value = "NYhrO8"

# oip column (U) changed from "null" to "-" for data rows 2-6
$ws.Range("U2:U6").Value = "-"

# hst column (V) changed from "0" to "-" for data rows 2-6, and its
# horizontal alignment now matches the left-aligned style used by the
# other text columns (e.g. U) instead of the previous right-aligned style.
$ws.Range("V2:V6").Value = "-"
$ws.Range("V2:V6").HorizontalAlignment = -4131
